$d = $word.ActiveDocument
$range1 = $d.Range(0, 4)
$range1.Text = "6.3 "

# source: the trailing space of run1 itself (position 3-4), which is color=2E54A5 spacing=12
$src = $d.Range(3, 4)
Write-Host ("src text=[" + $src.Text + "] len=" + $src.Text.Length)

$insPoint = $d.Range(4, 4)
$insPoint.FormattedText = $src
Write-Host ("Full: [" + $d.Content.Text.Substring(0,10) + "]")
